$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 <= original row 21
$ws.Range("B20").Value = 7188987
$ws.Range("E20").Value = 'NK Bistra'
$ws.Range("F20").Value = 'NK Zagorec Krapina'
$ws.Range("H20").Value = 1
$ws.Range("J20").Value = 3.4
$ws.Range("K20").Value = 3.4
$ws.Range("L20").Value = 1.909
$ws.Range("M20").Value = 3.5
$ws.Range("N20").Value = 3.4
$ws.Range("O20").Value = 1.909
$ws.Range("P20").Value = 0.5
$ws.Range("Q20").Value = 1.85
$ws.Range("R20").Value = 1.95
$ws.Range("S20").Value = 3
$ws.Range("X20").Value = 0.909
$ws.Range("Z20").Value = 0.95

# Row 21 <= original row 20
$ws.Range("B21").Value = 7188990
$ws.Range("E21").Value = 'NK Granicar Zupanja'
$ws.Range("F21").Value = 'NK Bedem Ivankovo'
$ws.Range("H21").Value = 2
$ws.Range("J21").Value = 2.1
$ws.Range("K21").Value = 3.5
$ws.Range("L21").Value = 3
$ws.Range("M21").Value = 2.1
$ws.Range("N21").Value = 3.5
$ws.Range("O21").Value = 3
$ws.Range("P21").Value = -0.25
$ws.Range("Q21").Value = 1.875
$ws.Range("R21").Value = 1.925
$ws.Range("S21").Value = 2.75
$ws.Range("X21").Value = 2
$ws.Range("Z21").Value = 0.925

# Row 65 <= original row 67
$ws.Range("B65").Value = 7493772
$ws.Range("E65").Value = 'NK Bedem Ivankovo'
$ws.Range("F65").Value = 'Sloga Nova Gradiska'
$ws.Range("H65").Value = 1
$ws.Range("I65").Value = 'H'
$ws.Range("J65").Value = 1.909
$ws.Range("K65").Value = 3.6
$ws.Range("L65").Value = 3.2
$ws.Range("M65").Value = 1.4
$ws.Range("N65").Value = 4.2
$ws.Range("O65").Value = 6.5
$ws.Range("P65").Value = -1.25
$ws.Range("Q65").Value = 1.9
$ws.Range("R65").Value = 1.9
$ws.Range("S65").Value = 3
$ws.Range("T65").Value = 1.85
$ws.Range("U65").Value = 1.95
$ws.Range("V65").Value = 0.3999999999999999
$ws.Range("W65").Value = -1
$ws.Range("Y65").Value = -0.5
$ws.Range("Z65").Value = 0.45
$ws.Range("AA65").Value = 0
$ws.Range("AB65").Value = 0

# Row 67 <= original row 65
$ws.Range("B67").Value = 7493773
$ws.Range("E67").Value = 'NK Zadar'
$ws.Range("F67").Value = 'NK Vodice'
$ws.Range("H67").Value = 2
$ws.Range("I67").Value = 'D'
$ws.Range("J67").Value = 1.2
$ws.Range("K67").Value = 6
$ws.Range("L67").Value = 9
$ws.Range("M67").Value = 1.2
$ws.Range("N67").Value = 6
$ws.Range("O67").Value = 9
$ws.Range("P67").Value = -2
$ws.Range("Q67").Value = 1.85
$ws.Range("R67").Value = 1.95
$ws.Range("S67").Value = 3.75
$ws.Range("T67").Value = 1.95
$ws.Range("U67").Value = 1.85
$ws.Range("V67").Value = -1
$ws.Range("W67").Value = 5
$ws.Range("Y67").Value = -1
$ws.Range("Z67").Value = 0.95
$ws.Range("AA67").Value = 0.475
$ws.Range("AB67").Value = -0.5

# Row 126 <= original row 129
$ws.Range("B126").Value = 8163880
$ws.Range("E126").Value = 'NK Dinamo Odranski Obre'
$ws.Range("F126").Value = 'Sava Strmec'
$ws.Range("H126").Value = 2
$ws.Range("I126").Value = 'D'
$ws.Range("J126").Value = 1.909
$ws.Range("K126").Value = 3.4
$ws.Range("L126").Value = 3.4
$ws.Range("M126").Value = 1.8
$ws.Range("N126").Value = 3.6
$ws.Range("O126").Value = 3.75
$ws.Range("P126").Value = -0.5
$ws.Range("Q126").Value = 1.825
$ws.Range("R126").Value = 1.975
$ws.Range("S126").Value = 3.25
$ws.Range("T126").Value = 1.925
$ws.Range("U126").Value = 1.875
$ws.Range("V126").Value = -1
$ws.Range("W126").Value = 2.6
$ws.Range("Y126").Value = -1
$ws.Range("Z126").Value = 0.9750000000000001
$ws.Range("AA126").Value = 0.925
$ws.Range("AB126").Value = -1

# Row 128 <= original row 126
$ws.Range("B128").Value = 8163883
$ws.Range("E128").Value = 'NK Zelina'
$ws.Range("F128").Value = 'NK Mladost Petrinja'
$ws.Range("G128").Value = 2
$ws.Range("I128").Value = 'H'
$ws.Range("J128").Value = 2
$ws.Range("K128").Value = 3.3
$ws.Range("L128").Value = 3.25
$ws.Range("M128").Value = 2.05
$ws.Range("N128").Value = 3.5
$ws.Range("O128").Value = 3
$ws.Range("P128").Value = -0.25
$ws.Range("Q128").Value = 1.85
$ws.Range("R128").Value = 1.95
$ws.Range("S128").Value = 3.5
$ws.Range("T128").Value = 1.825
$ws.Range("U128").Value = 1.975
$ws.Range("V128").Value = 1.05
$ws.Range("W128").Value = -1
$ws.Range("Y128").Value = 0.8500000000000001
$ws.Range("Z128").Value = -1
$ws.Range("AB128").Value = 0.9750000000000001

# Row 129 <= original row 128
$ws.Range("B129").Value = 8163879
$ws.Range("E129").Value = 'NK Bistra'
$ws.Range("F129").Value = 'NK Tresnjevka'
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 1
$ws.Range("J129").Value = 2.5
$ws.Range("L129").Value = 2.4
$ws.Range("M129").Value = 1.909
$ws.Range("O129").Value = 3.3
$ws.Range("Q129").Value = 1.975
$ws.Range("R129").Value = 1.825
$ws.Range("S129").Value = 2.75
$ws.Range("T129").Value = 1.75
$ws.Range("U129").Value = 1.95
$ws.Range("Z129").Value = 0.825
$ws.Range("AA129").Value = -1
$ws.Range("AB129").Value = 0.95

# Row 141 <= original row 142
$ws.Range("B141").Value = 8233939
$ws.Range("E141").Value = 'Sava Strmec'
$ws.Range("F141").Value = 'NK Mladost Petrinja'
$ws.Range("G141").Value = 0
$ws.Range("I141").Value = 'D'
$ws.Range("J141").Value = 2.25
$ws.Range("L141").Value = 2.55
$ws.Range("M141").Value = 2.25
$ws.Range("N141").Value = 3.6
$ws.Range("O141").Value = 2.55
$ws.Range("P141").Value = 0
$ws.Range("Q141").Value = 1.775
$ws.Range("R141").Value = 2.025
$ws.Range("T141").Value = 1.825
$ws.Range("U141").Value = 1.975
$ws.Range("V141").Value = -1
$ws.Range("W141").Value = 2.6
$ws.Range("Y141").Value = 0
$ws.Range("Z141").Value = 0
$ws.Range("AB141").Value = 0.9750000000000001

# Row 142 <= original row 141
$ws.Range("B142").Value = 8233938
$ws.Range("E142").Value = 'NK Bistra'
$ws.Range("F142").Value = 'NK Udarnik Kurilovec'
$ws.Range("G142").Value = 1
$ws.Range("I142").Value = 'H'
$ws.Range("J142").Value = 1.909
$ws.Range("L142").Value = 3.2
$ws.Range("M142").Value = 1.8
$ws.Range("N142").Value = 3.75
$ws.Range("O142").Value = 3.75
$ws.Range("P142").Value = -0.5
$ws.Range("Q142").Value = 1.975
$ws.Range("R142").Value = 1.825
$ws.Range("T142").Value = 1.775
$ws.Range("U142").Value = 2.025
$ws.Range("V142").Value = 0.8
$ws.Range("W142").Value = -1
$ws.Range("Y142").Value = 0.9750000000000001
$ws.Range("Z142").Value = -1
$ws.Range("AB142").Value = 1.025

